$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.85%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.96%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.555"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.82%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08037"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.98%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.906"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.274"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.90%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-10.48%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9449"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.38%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-2.20%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-4.17%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09706"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.49%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.04357"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.86%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1067"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001271"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.38%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005970"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.74%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.405"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.73%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'BitpandaEcosystemToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'0.3490"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.31%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'MCDex"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'9.974"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'14.41%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'ProBitToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.1379"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.82%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'ZBToken"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'0.2508"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.41%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'CoinExToken"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.04197"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-4.21%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'BitKan"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.001246"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.56%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'HotbitToken"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.004281"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.63%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001260"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.16%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003991"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.23%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02655"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-5.34%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05489"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-4.25%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007583"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-4.19%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1396"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.59%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007987"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-18.24%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002009"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.79%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008834"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-11.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-5.37%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.24%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.002271"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.23%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.005850"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'69.99%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.24%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.24%"
$ws.Range("E50").Style = "Normal"
